$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: safe to assign directly.
$textValues = @{
    "D2" = "68.125.13"
    "E2" = "  +3.71%  "
    "D3" = "3.343.98"
    "E3" = "  +1.18%  "
    "E4" = "  -0.28%  "
    "E6" = "  -0.19%  "
    "E7" = "  +0.04%  "
    "E8" = "  +3.39%  "
    "E9" = "  +3.47%  "
    "E10" = "  +1.33%  "
    "E11" = "  +2.49%  "
    "E12" = "  +6.59%  "
    "E13" = "  +11.59%  "
    "D14" = "3.900.74"
    "E14" = "  +1.57%  "
    "E15" = "  +0.99%  "
    "D16" = "68.325.61"
    "E16" = "  +4.02%  "
    "E17" = "  +1.73%  "
    "D18" = "3.367.84"
    "E18" = "  +1.57%  "
    "E19" = "  +1.02%  "
    "E20" = "  +1.88%  "
    "E21" = "  +1.84%  "
    "E22" = "  -0.24%  "
    "E23" = "  +1.72%  "
    "E24" = "  +0.37%  "
    "E25" = "  +3.57%  "
    "E26" = "  +5.65%  "
    "E27" = "  +4.04%  "
    "E28" = "  +8.54%  "
    "E29" = "  +2.08%  "
    "E30" = "  +1.48%  "
    "E31" = "  +8.86%  "
    "D32" = "3.977.65"
    "E32" = "  +6.22%  "
    "E33" = "  -0.47%  "
    "E34" = "  +2.10%  "
    "E35" = "  +2.39%  "
    "E36" = "  -0.11%  "
    "E37" = "  +0.81%  "
    "E38" = "  +6.70%  "
    "B39" = "Kaspa"
    "C39" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "E39" = "  +4.12%  "
    "B40" = "Stacks"
    "C40" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "E40" = "  +5.54%  "
    "E41" = "  -0.51%  "
    "D42" = "0.0₃0698"
    "E42" = "  +2.74%  "
    "E43" = "  +3.04%  "
    "E44" = "  +2.97%  "
    "E45" = "  +3.02%  "
    "E46" = "  +2.36%  "
    "B47" = "FirstDigitalUSD"
    "C47" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "E47" = "  +0.39%  "
    "B48" = "ThetaToken"
    "C48" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "E48" = "  +2.43%  "
    "E49" = "  +9.21%  "
    "E50" = "  +5.22%  "
    "E51" = "  +6.83%  "
}

# Numeric-looking Price values that must stay TEXT (preserve trailing zeros / exact formatting).
$numericTextValues = @{
    "D4" = "1.00"
    "D5" = "592.35"
    "D6" = "184.65"
    "D8" = "0.593"
    "D9" = "0.180"
    "D10" = "0.583"
    "D11" = "46.82"
    "D12" = "0.0000277"
    "D13" = "639.45"
    "D15" = "8.47"
    "D19" = "17.78"
    "D20" = "10.99"
    "D21" = "0.904"
    "D22" = "17.86"
    "D23" = "5.05"
    "D24" = "98.35"
    "D25" = "4.07"
    "D26" = "2.82"
    "D27" = "9.68"
    "D29" = "8.61"
    "D30" = "6.75"
    "D31" = "605.80"
    "D33" = "3.65"
    "D34" = "11.01"
    "D36" = "0.998"
    "D37" = "55.96"
    "D38" = "2.74"
    "D39" = "0.130"
    "D40" = "3.28"
    "D41" = "33.28"
    "D44" = "0.340"
    "D45" = "0.0419"
    "D47" = "1.01"
    "D48" = "2.56"
    "D50" = "131.34"
    "D51" = "7.72"
}

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}

foreach ($addr in $numericTextValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextValues[$addr]
    $cell.Style = "Normal"
}
